# Add two new sub-bullets under the "obstacle avoidance" list item:
#   - "Seemingly is the opposite to the centre track."      (ilvl = 2)
#   - "IF obstacle is to the right move to the left"        (ilvl = 3)
# They are inserted right after the existing bullet
#   "Add rules to the engine such as in the c++ example for it to avoid obstacles."
# and right before "Make the Current multiple movement engines into the one engine!".

$d = $word.ActiveDocument

# Locate the anchor paragraph robustly via Find rather than a hard-coded index.
$anchorText = "Add rules to the engine such as in the c++ example for it to avoid obstacles."
$searchRange = $d.Content
$found = $searchRange.Find.Execute($anchorText, $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find anchor paragraph: $anchorText"
}

$anchorPara = $searchRange.Paragraphs(1)

# Insert the first new paragraph right after the anchor paragraph.
$null = $anchorPara.Range.InsertParagraphAfter()
$p1 = $anchorPara.Next()
$p1.Range.Text = "Seemingly is the opposite to the centre track."
$p1.Range.ListFormat.ListLevelNumber = 3   # ListLevelNumber is 1-based -> w:ilvl="2"

# Insert the second new paragraph right after the first new one.
$null = $p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "IF obstacle is to the right move to the left"
$p2.Range.ListFormat.ListLevelNumber = 4   # ListLevelNumber is 1-based -> w:ilvl="3"

Write-Host "Inserted: '$($p1.Range.Text)' (ilvl=$($p1.Range.ListFormat.ListLevelNumber - 1))"
Write-Host "Inserted: '$($p2.Range.Text)' (ilvl=$($p2.Range.ListFormat.ListLevelNumber - 1))"
